$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (74 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 443  # H17: 1937.6364 -> 443
$ws.Cells.Item(17, 10).Value = 439.27118  # J17: 1957.2461 -> 439.27118
$ws.Cells.Item(17, 12).Value = 1317.81354  # L17: 5871.7383 -> 1317.81354
$ws.Cells.Item(17, 14).Value = -1653.81354  # N17: -6207.7383 -> -1653.81354
$ws.Cells.Item(33, 8).Value = 422.4375  # H33: 437.33334 -> 422.4375
$ws.Cells.Item(33, 9).Value = 368.5  # I33: 381.53845 -> 368.5
$ws.Cells.Item(33, 11).Value = 368.5  # K33: 381.53845 -> 368.5
$ws.Cells.Item(33, 13).Value = -139.5  # M33: -152.53845 -> -139.5
$ws.Cells.Item(43, 8).Value = 700.6923  # H43: 686.2857 -> 700.6923
$ws.Cells.Item(43, 10).Value = 719.8889  # J43: 697.8 -> 719.8889
$ws.Cells.Item(43, 12).Value = 719.8889  # L43: 697.8 -> 719.8889
$ws.Cells.Item(43, 14).Value = -857.8889  # N43: -835.8 -> -857.8889
$ws.Cells.Item(62, 8).Value = 2355.4707  # H62: 2403 -> 2355.4707
$ws.Cells.Item(62, 10).Value = 2493.125  # J62: 2657.8333 -> 2493.125
$ws.Cells.Item(62, 12).Value = 2493.125  # L62: 2657.8333 -> 2493.125
$ws.Cells.Item(62, 14).Value = -3741.125  # N62: -3905.8333 -> -3741.125
$ws.Cells.Item(65, 8).Value = 2355.4707  # H65: 2403 -> 2355.4707
$ws.Cells.Item(65, 10).Value = 2493.125  # J65: 2657.8333 -> 2493.125
$ws.Cells.Item(65, 12).Value = 12465.625  # L65: 13289.1665 -> 12465.625
$ws.Cells.Item(65, 14).Value = -18705.625  # N65: -19529.1665 -> -18705.625
$ws.Cells.Item(70, 8).Value = 9645.272000000001  # H70: 8949.833000000001 -> 9645.272000000001
$ws.Cells.Item(70, 9).Value = 20440  # I70: 14828.571 -> 20440
$ws.Cells.Item(70, 10).Value = 649.6667  # J70: 719.6 -> 649.6667
$ws.Cells.Item(70, 11).Value = 61320  # K70: 44485.713 -> 61320
$ws.Cells.Item(70, 12).Value = 1949.0001  # L70: 2158.8 -> 1949.0001
$ws.Cells.Item(70, 13).Value = -61050  # M70: -44215.713 -> -61050
$ws.Cells.Item(70, 14).Value = -2489.0001  # N70: -2698.8 -> -2489.0001
$ws.Cells.Item(73, 8).Value = 9645.272000000001  # H73: 8949.833000000001 -> 9645.272000000001
$ws.Cells.Item(73, 9).Value = 20440  # I73: 14828.571 -> 20440
$ws.Cells.Item(73, 10).Value = 649.6667  # J73: 719.6 -> 649.6667
$ws.Cells.Item(73, 11).Value = 61320  # K73: 44485.713 -> 61320
$ws.Cells.Item(73, 12).Value = 1949.0001  # L73: 2158.8 -> 1949.0001
$ws.Cells.Item(73, 13).Value = -60384  # M73: -43549.713 -> -60384
$ws.Cells.Item(73, 14).Value = -3821.0001  # N73: -4030.8 -> -3821.0001
$ws.Cells.Item(88, 8).Value = 1935.6666  # H88: 2002.6666 -> 1935.6666
$ws.Cells.Item(88, 9).Value = 1901.5  # I88: 2000 -> 1901.5
$ws.Cells.Item(88, 11).Value = 1901.5  # K88: 2000 -> 1901.5
$ws.Cells.Item(88, 13).Value = -1495.5  # M88: -1594 -> -1495.5
$ws.Cells.Item(91, 8).Value = 1935.6666  # H91: 2002.6666 -> 1935.6666
$ws.Cells.Item(91, 9).Value = 1901.5  # I91: 2000 -> 1901.5
$ws.Cells.Item(91, 11).Value = 1901.5  # K91: 2000 -> 1901.5
$ws.Cells.Item(91, 13).Value = -497.5  # M91: -596 -> -497.5
$ws.Cells.Item(115, 8).Value = 692.5  # H115: 649 -> 692.5
$ws.Cells.Item(115, 9).Value = 692.5  # I115: 649 -> 692.5
$ws.Cells.Item(115, 11).Value = 2077.5  # K115: 1947 -> 2077.5
$ws.Cells.Item(115, 13).Value = -510.5  # M115: -380 -> -510.5
$ws.Cells.Item(116, 8).Value = 4530.154  # H116: 3823.1667 -> 4530.154
$ws.Cells.Item(116, 9).Value = 1693.3334  # I116: 1750.8334 -> 1693.3334
$ws.Cells.Item(116, 10).Value = 5381.2  # J116: 4859.3335 -> 5381.2
$ws.Cells.Item(116, 11).Value = 1693.3334  # K116: 1750.8334 -> 1693.3334
$ws.Cells.Item(116, 12).Value = 5381.2  # L116: 4859.3335 -> 5381.2
$ws.Cells.Item(116, 13).Value = 1748.6666  # M116: 1691.1666 -> 1748.6666
$ws.Cells.Item(116, 14).Value = -12265.2  # N116: -11743.3335 -> -12265.2
$ws.Cells.Item(118, 8).Value = 745  # H118: 700 -> 745
$ws.Cells.Item(118, 9).Value = 745  # I118: 0 -> 745
$ws.Cells.Item(118, 10).Value = 0  # J118: 700 -> 0
$ws.Cells.Item(118, 11).Value = 2235  # K118: 0 -> 2235
$ws.Cells.Item(118, 12).Value = 0  # L118: 2100 -> 0
$ws.Cells.Item(118, 13).Value = -578  # M118: None -> -578
$ws.Cells.Item(118, 14).Value = $null  # N118: -5414 -> None
$ws.Cells.Item(132, 8).Value = 34931.453  # H132: 29276.621 -> 34931.453
$ws.Cells.Item(132, 9).Value = 41383.42  # I132: 32606.94 -> 41383.42
$ws.Cells.Item(132, 10).Value = 1381.2  # J132: 1801.5 -> 1381.2
$ws.Cells.Item(132, 11).Value = 124150.26  # K132: 97820.81999999999 -> 124150.26
$ws.Cells.Item(132, 12).Value = 4143.6  # L132: 5404.5 -> 4143.6
$ws.Cells.Item(132, 13).Value = -121620.26  # M132: -95290.81999999999 -> -121620.26
$ws.Cells.Item(132, 14).Value = -9203.6  # N132: -10464.5 -> -9203.6
$ws.Cells.Item(138, 8).Value = 2753.9111  # H138: 2788.9575 -> 2753.9111
$ws.Cells.Item(138, 9).Value = 1224.1  # I138: 1300.125 -> 1224.1
$ws.Cells.Item(138, 10).Value = 3191  # J138: 3094.359 -> 3191
$ws.Cells.Item(138, 11).Value = 3672.3  # K138: 3900.375 -> 3672.3
$ws.Cells.Item(138, 12).Value = 9573  # L138: 9283.076999999999 -> 9573
$ws.Cells.Item(138, 13).Value = 1467.7  # M138: 1239.625 -> 1467.7
$ws.Cells.Item(138, 14).Value = -19853  # N138: -19563.077 -> -19853

# --- Sheet: ARM (29 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17194.5  # H32: 20131.81 -> 17194.5
$ws.Cells.Item(32, 9).Value = 20302  # I32: 24880.158 -> 20302
$ws.Cells.Item(32, 11).Value = 20302  # K32: 24880.158 -> 20302
$ws.Cells.Item(32, 13).Value = -20015  # M32: -24593.158 -> -20015
$ws.Cells.Item(45, 8).Value = 2654.195  # H45: 2642.3171 -> 2654.195
$ws.Cells.Item(45, 9).Value = 2066.1875  # I45: 1963 -> 2066.1875
$ws.Cells.Item(45, 10).Value = 3030.52  # J45: 3123.5 -> 3030.52
$ws.Cells.Item(45, 11).Value = 2066.1875  # K45: 1963 -> 2066.1875
$ws.Cells.Item(45, 12).Value = 3030.52  # L45: 3123.5 -> 3030.52
$ws.Cells.Item(45, 13).Value = -1689.1875  # M45: -1586 -> -1689.1875
$ws.Cells.Item(45, 14).Value = -3784.52  # N45: -3877.5 -> -3784.52
$ws.Cells.Item(122, 8).Value = 2135.1667  # H122: 2536.7585 -> 2135.1667
$ws.Cells.Item(122, 9).Value = 2131.2693  # I122: 2543.476 -> 2131.2693
$ws.Cells.Item(122, 10).Value = 2145.3  # J122: 2519.125 -> 2145.3
$ws.Cells.Item(122, 11).Value = 6393.8079  # K122: 7630.428 -> 6393.8079
$ws.Cells.Item(122, 12).Value = 6435.900000000001  # L122: 7557.375 -> 6435.900000000001
$ws.Cells.Item(122, 13).Value = -3943.8079  # M122: -5180.428 -> -3943.8079
$ws.Cells.Item(122, 14).Value = -11335.9  # N122: -12457.375 -> -11335.9
$ws.Cells.Item(132, 8).Value = 14302.902  # H132: 15069.41 -> 14302.902
$ws.Cells.Item(132, 9).Value = 2120.6  # I132: 2253.4062 -> 2120.6
$ws.Cells.Item(132, 10).Value = 85366.336  # J132: 73656.86 -> 85366.336
$ws.Cells.Item(132, 11).Value = 6361.799999999999  # K132: 6760.2186 -> 6361.799999999999
$ws.Cells.Item(132, 12).Value = 256099.008  # L132: 220970.58 -> 256099.008
$ws.Cells.Item(132, 13).Value = -3831.799999999999  # M132: -4230.2186 -> -3831.799999999999
$ws.Cells.Item(132, 14).Value = -261159.008  # N132: -226030.58 -> -261159.008
$ws.Cells.Item(138, 8).Value = 53429  # H138: 50329 -> 53429
$ws.Cells.Item(138, 10).Value = 53429  # J138: 50329 -> 53429
$ws.Cells.Item(138, 12).Value = 53429  # L138: 50329 -> 53429
$ws.Cells.Item(138, 14).Value = -63709  # N138: -60609 -> -63709

# --- Sheet: BSM (40 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 19841.715  # H81: 15579.6 -> 19841.715
$ws.Cells.Item(81, 10).Value = 19841.715  # J81: 15579.6 -> 19841.715
$ws.Cells.Item(81, 12).Value = 19841.715  # L81: 15579.6 -> 19841.715
$ws.Cells.Item(81, 14).Value = -21963.715  # N81: -17701.6 -> -21963.715
$ws.Cells.Item(84, 8).Value = 19841.715  # H84: 15579.6 -> 19841.715
$ws.Cells.Item(84, 10).Value = 19841.715  # J84: 15579.6 -> 19841.715
$ws.Cells.Item(84, 12).Value = 59525.145  # L84: 46738.8 -> 59525.145
$ws.Cells.Item(84, 14).Value = -70133.145  # N84: -57346.8 -> -70133.145
$ws.Cells.Item(86, 8).Value = 1722.2222  # H86: 1837.5 -> 1722.2222
$ws.Cells.Item(86, 9).Value = 1540  # I86: 1630 -> 1540
$ws.Cells.Item(86, 10).Value = 2633.3333  # J86: 2183.3333 -> 2633.3333
$ws.Cells.Item(86, 11).Value = 1540  # K86: 1630 -> 1540
$ws.Cells.Item(86, 12).Value = 2633.3333  # L86: 2183.3333 -> 2633.3333
$ws.Cells.Item(86, 13).Value = -417  # M86: -507 -> -417
$ws.Cells.Item(86, 14).Value = -4879.3333  # N86: -4429.3333 -> -4879.3333
$ws.Cells.Item(89, 8).Value = 1722.2222  # H89: 1837.5 -> 1722.2222
$ws.Cells.Item(89, 9).Value = 1540  # I89: 1630 -> 1540
$ws.Cells.Item(89, 10).Value = 2633.3333  # J89: 2183.3333 -> 2633.3333
$ws.Cells.Item(89, 11).Value = 7700  # K89: 8150 -> 7700
$ws.Cells.Item(89, 12).Value = 13166.6665  # L89: 10916.6665 -> 13166.6665
$ws.Cells.Item(89, 13).Value = -2084  # M89: -2534 -> -2084
$ws.Cells.Item(89, 14).Value = -24398.6665  # N89: -22148.6665 -> -24398.6665
$ws.Cells.Item(94, 8).Value = 3673.25  # H94: 3961.2727 -> 3673.25
$ws.Cells.Item(94, 9).Value = 1684.2222  # I94: 1868.5 -> 1684.2222
$ws.Cells.Item(94, 10).Value = 4866.6665  # J94: 5157.143 -> 4866.6665
$ws.Cells.Item(94, 11).Value = 1684.2222  # K94: 1868.5 -> 1684.2222
$ws.Cells.Item(94, 12).Value = 4866.6665  # L94: 5157.143 -> 4866.6665
$ws.Cells.Item(94, 13).Value = -1233.2222  # M94: -1417.5 -> -1233.2222
$ws.Cells.Item(94, 14).Value = -5768.6665  # N94: -6059.143 -> -5768.6665
$ws.Cells.Item(112, 8).Value = 44821.668  # H112: 50000 -> 44821.668
$ws.Cells.Item(112, 10).Value = 44821.668  # J112: 50000 -> 44821.668
$ws.Cells.Item(112, 12).Value = 44821.668  # L112: 50000 -> 44821.668
$ws.Cells.Item(112, 14).Value = -47775.668  # N112: -52954 -> -47775.668
$ws.Cells.Item(134, 8).Value = 59740.223  # H134: 118691.555 -> 59740.223
$ws.Cells.Item(134, 9).Value = 63207.293  # I134: 118691.555 -> 63207.293
$ws.Cells.Item(134, 10).Value = 800  # J134: 0 -> 800
$ws.Cells.Item(134, 11).Value = 189621.879  # K134: 356074.665 -> 189621.879
$ws.Cells.Item(134, 12).Value = 2400  # L134: 0 -> 2400
$ws.Cells.Item(134, 13).Value = -187086.879  # M134: -353539.665 -> -187086.879
$ws.Cells.Item(134, 14).Value = -7470  # N134: None -> -7470

# --- Sheet: CRP (21 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9555.089  # H31: 10090.038 -> 9555.089
$ws.Cells.Item(31, 9).Value = 13304.471  # I31: 15034.667 -> 13304.471
$ws.Cells.Item(31, 10).Value = 3760.5908  # J31: 3640.5217 -> 3760.5908
$ws.Cells.Item(31, 11).Value = 13304.471  # K31: 15034.667 -> 13304.471
$ws.Cells.Item(31, 12).Value = 3760.5908  # L31: 3640.5217 -> 3760.5908
$ws.Cells.Item(31, 13).Value = -13009.471  # M31: -14739.667 -> -13009.471
$ws.Cells.Item(31, 14).Value = -4350.5908  # N31: -4230.521699999999 -> -4350.5908
$ws.Cells.Item(34, 8).Value = 9555.089  # H34: 10090.038 -> 9555.089
$ws.Cells.Item(34, 9).Value = 13304.471  # I34: 15034.667 -> 13304.471
$ws.Cells.Item(34, 10).Value = 3760.5908  # J34: 3640.5217 -> 3760.5908
$ws.Cells.Item(34, 11).Value = 13304.471  # K34: 15034.667 -> 13304.471
$ws.Cells.Item(34, 12).Value = 3760.5908  # L34: 3640.5217 -> 3760.5908
$ws.Cells.Item(34, 13).Value = -13102.471  # M34: -14832.667 -> -13102.471
$ws.Cells.Item(34, 14).Value = -4164.5908  # N34: -4044.5217 -> -4164.5908
$ws.Cells.Item(134, 8).Value = 6691  # H134: 7673.6665 -> 6691
$ws.Cells.Item(134, 9).Value = 716.4666999999999  # I134: 842.0833 -> 716.4666999999999
$ws.Cells.Item(134, 10).Value = 51500  # J134: 35000 -> 51500
$ws.Cells.Item(134, 11).Value = 2149.4001  # K134: 2526.2499 -> 2149.4001
$ws.Cells.Item(134, 12).Value = 154500  # L134: 105000 -> 154500
$ws.Cells.Item(134, 13).Value = 385.5999000000002  # M134: 8.750100000000202 -> 385.5999000000002
$ws.Cells.Item(134, 14).Value = -159570  # N134: -110070 -> -159570

# --- Sheet: CUL (80 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 825.55554  # H34: 850 -> 825.55554
$ws.Cells.Item(34, 9).Value = 180  # I34: 1000 -> 180
$ws.Cells.Item(34, 10).Value = 906.25  # J34: 812.5 -> 906.25
$ws.Cells.Item(34, 11).Value = 540  # K34: 3000 -> 540
$ws.Cells.Item(34, 12).Value = 2718.75  # L34: 2437.5 -> 2718.75
$ws.Cells.Item(34, 13).Value = -456  # M34: -2916 -> -456
$ws.Cells.Item(34, 14).Value = -2886.75  # N34: -2605.5 -> -2886.75
$ws.Cells.Item(68, 8).Value = 1300.25  # H68: 1234.3715 -> 1300.25
$ws.Cells.Item(68, 9).Value = 1500  # I68: 749.75 -> 1500
$ws.Cells.Item(68, 10).Value = 1293.8064  # J68: 1296.9032 -> 1293.8064
$ws.Cells.Item(68, 11).Value = 4500  # K68: 2249.25 -> 4500
$ws.Cells.Item(68, 12).Value = 3881.4192  # L68: 3890.7096 -> 3881.4192
$ws.Cells.Item(68, 13).Value = -3689  # M68: -1438.25 -> -3689
$ws.Cells.Item(68, 14).Value = -5503.4192  # N68: -5512.7096 -> -5503.4192
$ws.Cells.Item(70, 8).Value = 3531  # H70: 4351.8184 -> 3531
$ws.Cells.Item(70, 9).Value = 2641.2  # I70: 3800 -> 2641.2
$ws.Cells.Item(70, 11).Value = 7923.599999999999  # K70: 11400 -> 7923.599999999999
$ws.Cells.Item(70, 13).Value = -7608.599999999999  # M70: -11085 -> -7608.599999999999
$ws.Cells.Item(71, 8).Value = 1300.25  # H71: 1234.3715 -> 1300.25
$ws.Cells.Item(71, 9).Value = 1500  # I71: 749.75 -> 1500
$ws.Cells.Item(71, 10).Value = 1293.8064  # J71: 1296.9032 -> 1293.8064
$ws.Cells.Item(71, 11).Value = 13500  # K71: 6747.75 -> 13500
$ws.Cells.Item(71, 12).Value = 11644.2576  # L71: 11672.1288 -> 11644.2576
$ws.Cells.Item(71, 13).Value = -9444  # M71: -2691.75 -> -9444
$ws.Cells.Item(71, 14).Value = -19756.2576  # N71: -19784.1288 -> -19756.2576
$ws.Cells.Item(73, 8).Value = 3531  # H73: 4351.8184 -> 3531
$ws.Cells.Item(73, 9).Value = 2641.2  # I73: 3800 -> 2641.2
$ws.Cells.Item(73, 11).Value = 7923.599999999999  # K73: 11400 -> 7923.599999999999
$ws.Cells.Item(73, 13).Value = -6831.599999999999  # M73: -10308 -> -6831.599999999999
$ws.Cells.Item(75, 8).Value = 983  # H75: 1953.75 -> 983
$ws.Cells.Item(75, 10).Value = 983  # J75: 1953.75 -> 983
$ws.Cells.Item(75, 12).Value = 2949  # L75: 5861.25 -> 2949
$ws.Cells.Item(75, 14).Value = -4945  # N75: -7857.25 -> -4945
$ws.Cells.Item(78, 8).Value = 983  # H78: 1953.75 -> 983
$ws.Cells.Item(78, 10).Value = 983  # J78: 1953.75 -> 983
$ws.Cells.Item(78, 12).Value = 8847  # L78: 17583.75 -> 8847
$ws.Cells.Item(78, 14).Value = -18831  # N78: -27567.75 -> -18831
$ws.Cells.Item(103, 8).Value = 930.2  # H103: 2178.4285 -> 930.2
$ws.Cells.Item(103, 9).Value = 384.8  # I103: 630 -> 384.8
$ws.Cells.Item(103, 10).Value = 1475.6  # J103: 6049.5 -> 1475.6
$ws.Cells.Item(103, 11).Value = 1154.4  # K103: 1890 -> 1154.4
$ws.Cells.Item(103, 12).Value = 4426.799999999999  # L103: 18148.5 -> 4426.799999999999
$ws.Cells.Item(103, 13).Value = -275.4000000000001  # M103: -1011 -> -275.4000000000001
$ws.Cells.Item(103, 14).Value = -6184.799999999999  # N103: -19906.5 -> -6184.799999999999
$ws.Cells.Item(107, 8).Value = 10245.272  # H107: 9053.462 -> 10245.272
$ws.Cells.Item(107, 9).Value = 25300  # I107: 33600 -> 25300
$ws.Cells.Item(107, 10).Value = 1642.5714  # J107: 1689.5 -> 1642.5714
$ws.Cells.Item(107, 11).Value = 75900  # K107: 100800 -> 75900
$ws.Cells.Item(107, 12).Value = 4927.7142  # L107: 5068.5 -> 4927.7142
$ws.Cells.Item(107, 13).Value = -73980  # M107: -98880 -> -73980
$ws.Cells.Item(107, 14).Value = -8767.7142  # N107: -8908.5 -> -8767.7142
$ws.Cells.Item(114, 8).Value = 1503.75  # H114: 1586.875 -> 1503.75
$ws.Cells.Item(114, 9).Value = 2412.5  # I114: 2099 -> 2412.5
$ws.Cells.Item(114, 10).Value = 595  # J114: 733.3333 -> 595
$ws.Cells.Item(114, 11).Value = 7237.5  # K114: 6297 -> 7237.5
$ws.Cells.Item(114, 12).Value = 1785  # L114: 2199.9999 -> 1785
$ws.Cells.Item(114, 13).Value = -3983.5  # M114: -3043 -> -3983.5
$ws.Cells.Item(114, 14).Value = -8293  # N114: -8707.999899999999 -> -8293
$ws.Cells.Item(117, 8).Value = 37038456  # H117: 37038880 -> 37038456
$ws.Cells.Item(117, 9).Value = 769.6  # I117: 829.6 -> 769.6
$ws.Cells.Item(117, 10).Value = 83335570  # J117: 83336450 -> 83335570
$ws.Cells.Item(117, 11).Value = 2308.8  # K117: 2488.8 -> 2308.8
$ws.Cells.Item(117, 12).Value = 250006710  # L117: 250009350 -> 250006710
$ws.Cells.Item(117, 13).Value = 1133.2  # M117: 953.1999999999998 -> 1133.2
$ws.Cells.Item(117, 14).Value = -250013594  # N117: -250016234 -> -250013594
$ws.Cells.Item(121, 8).Value = 3571.257  # H121: 3979.742 -> 3571.257
$ws.Cells.Item(121, 10).Value = 4369.7856  # J121: 5030.5 -> 4369.7856
$ws.Cells.Item(121, 12).Value = 13109.3568  # L121: 15091.5 -> 13109.3568
$ws.Cells.Item(121, 14).Value = -15729.3568  # N121: -17711.5 -> -15729.3568
$ws.Cells.Item(129, 8).Value = 358528.66  # H129: 278588.38 -> 358528.66
$ws.Cells.Item(129, 9).Value = 513.6  # I129: 765.4286 -> 513.6
$ws.Cells.Item(129, 10).Value = 557425.9  # J129: 455384.8 -> 557425.9
$ws.Cells.Item(129, 11).Value = 1540.8  # K129: 2296.2858 -> 1540.8
$ws.Cells.Item(129, 12).Value = 1672277.7  # L129: 1366154.4 -> 1672277.7
$ws.Cells.Item(129, 13).Value = 3459.2  # M129: 2703.7142 -> 3459.2
$ws.Cells.Item(129, 14).Value = -1682277.7  # N129: -1376154.4 -> -1682277.7
$ws.Cells.Item(131, 8).Value = 807.9400000000001  # H131: 805.39 -> 807.9400000000001
$ws.Cells.Item(131, 10).Value = 821.5  # J131: 818.84375 -> 821.5
$ws.Cells.Item(131, 12).Value = 2464.5  # L131: 2456.53125 -> 2464.5
$ws.Cells.Item(131, 14).Value = -12544.5  # N131: -12536.53125 -> -12544.5

# --- Sheet: GSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8626.700000000001  # H80: 9379.666999999999 -> 8626.700000000001
$ws.Cells.Item(80, 9).Value = 12625.5  # I80: 15319.375 -> 12625.5
$ws.Cells.Item(80, 11).Value = 12625.5  # K80: 15319.375 -> 12625.5
$ws.Cells.Item(80, 13).Value = -11627.5  # M80: -14321.375 -> -11627.5
$ws.Cells.Item(83, 8).Value = 8626.700000000001  # H83: 9379.666999999999 -> 8626.700000000001
$ws.Cells.Item(83, 9).Value = 12625.5  # I83: 15319.375 -> 12625.5
$ws.Cells.Item(83, 11).Value = 63127.5  # K83: 76596.875 -> 63127.5
$ws.Cells.Item(83, 13).Value = -58135.5  # M83: -71604.875 -> -58135.5
$ws.Cells.Item(126, 8).Value = 5505.4  # H126: 5233.1875 -> 5505.4
$ws.Cells.Item(126, 9).Value = 5136.778  # I126: 4723.1 -> 5136.778
$ws.Cells.Item(126, 10).Value = 6058.3335  # J126: 6083.3335 -> 6058.3335
$ws.Cells.Item(126, 11).Value = 15410.334  # K126: 14169.3 -> 15410.334
$ws.Cells.Item(126, 12).Value = 18175.0005  # L126: 18250.0005 -> 18175.0005
$ws.Cells.Item(126, 13).Value = -12940.334  # M126: -11699.3 -> -12940.334
$ws.Cells.Item(126, 14).Value = -23115.0005  # N126: -23190.0005 -> -23115.0005

# --- Sheet: LTW (14 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2038.7059  # H82: 2335.7144 -> 2038.7059
$ws.Cells.Item(82, 9).Value = 2336.5386  # I82: 3011.111 -> 2336.5386
$ws.Cells.Item(82, 10).Value = 1070.75  # J82: 1120 -> 1070.75
$ws.Cells.Item(82, 11).Value = 2336.5386  # K82: 3011.111 -> 2336.5386
$ws.Cells.Item(82, 12).Value = 1070.75  # L82: 1120 -> 1070.75
$ws.Cells.Item(82, 13).Value = -1975.5386  # M82: -2650.111 -> -1975.5386
$ws.Cells.Item(82, 14).Value = -1792.75  # N82: -1842 -> -1792.75
$ws.Cells.Item(85, 8).Value = 2038.7059  # H85: 2335.7144 -> 2038.7059
$ws.Cells.Item(85, 9).Value = 2336.5386  # I85: 3011.111 -> 2336.5386
$ws.Cells.Item(85, 10).Value = 1070.75  # J85: 1120 -> 1070.75
$ws.Cells.Item(85, 11).Value = 2336.5386  # K85: 3011.111 -> 2336.5386
$ws.Cells.Item(85, 12).Value = 1070.75  # L85: 1120 -> 1070.75
$ws.Cells.Item(85, 13).Value = -1088.5386  # M85: -1763.111 -> -1088.5386
$ws.Cells.Item(85, 14).Value = -3566.75  # N85: -3616 -> -3566.75
